# Update TPM-derived values in the LR-pairs sheet (Fga-Itga2b)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.2303363333333333
$ws.Range("N2").Value = 0.691009
$ws.Range("O2").Value = 0.0420565315194687
$ws.Range("P2").Value = 0.0420565315194687
$ws.Range("Q2").Value = 0.3055554270193334
$ws.Range("R2").Value = 2.749998843174
$ws.Range("S2").Value = 0.0420565315194687
$ws.Range("T2").Value = 0.0420565315194687

# Row 3
$ws.Range("O3").Value = 0.8440851393264226
$ws.Range("P3").Value = 0.8440851393264227
$ws.Range("S3").Value = 0.8440851393264226
$ws.Range("T3").Value = 0.8440851393264227

# Row 4
$ws.Range("M4").Value = 0.6235823333333333
$ws.Range("N4").Value = 1.870747
$ws.Range("O4").Value = 0.1138583291541087
$ws.Range("P4").Value = 0.1138583291541087
$ws.Range("Q4").Value = 0.8272206272713334
$ws.Range("R4").Value = 7.444985645442001
$ws.Range("S4").Value = 0.1138583291541087
$ws.Range("T4").Value = 0.1138583291541087
